# Apply updated cryptocurrency market data (prices / 1h volume change)
# as produced by the scheduled GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.634.63'
$ws.Range("E2").Value = '  -1.80%  '

$ws.Range("D3").Value = '3.002.00'
$ws.Range("E3").Value = '  -0.98%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '598.67'
$ws.Range("E5").Value = '  +2.62%  '

$ws.Range("D6").Value = '144.46'
$ws.Range("E6").Value = '  -3.38%  '

$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").Value = '0.522'
$ws.Range("E8").Value = '  -0.68%  '

$ws.Range("D9").Value = '3.002.80'
$ws.Range("E9").Value = '  -0.95%  '

$ws.Range("E10").Value = '  -2.79%  '

$ws.Range("E11").Value = '  +4.26%  '

$ws.Range("D12").Value = '0.464'
$ws.Range("E12").Value = '  +4.34%  '

$ws.Range("E13").Value = '  -1.30%  '

$ws.Range("D14").Value = '34.34'
$ws.Range("E14").Value = '  -2.99%  '

$ws.Range("E15").Value = '  +2.67%  '

$ws.Range("D16").Value = '3.491.55'
$ws.Range("E16").Value = '  -1.20%  '

$ws.Range("D17").Value = '7.02'
$ws.Range("E17").Value = '  -0.70%  '

$ws.Range("D18").Value = '61.540.51'
$ws.Range("E18").Value = '  -1.87%  '

$ws.Range("D19").Value = '2.999.17'
$ws.Range("E19").Value = '  -1.03%  '

$ws.Range("D20").Value = '455.34'
$ws.Range("E20").Value = '  -2.74%  '

$ws.Range("E21").Value = '  -0.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.690'
$ws.Range("E22").Value = '  -0.47%  '

$ws.Range("D23").Value = '7.37'
$ws.Range("E23").Value = '  -0.78%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.50'
$ws.Range("E24").Value = '  +1.72%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.20'
$ws.Range("E25").Value = '  -8.35%  '

$ws.Range("D26").Value = '12.21'
$ws.Range("E26").Value = '  -1.82%  '

$ws.Range("D27").Value = '10.44'
$ws.Range("E27").Value = '  -0.74%  '

$ws.Range("E28").Value = '  +0.08%  '

$ws.Range("D29").Value = '2.68'
$ws.Range("E29").Value = '  +1.86%  '

$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.08%  '

$ws.Range("E31").Value = '  -3.02%  '

$ws.Range("D32").Value = '2.06'
$ws.Range("E32").Value = '  -4.38%  '

$ws.Range("D33").Value = '27.32'
$ws.Range("E33").Value = '  -0.88%  '

$ws.Range("E34").Value = '  -0.30%  '

$ws.Range("D35").Value = '0.0₃0819'
$ws.Range("E35").Value = '  +2.60%  '

$ws.Range("E36").Value = '  -1.97%  '

$ws.Range("D37").Value = '5.76'
$ws.Range("E37").Value = '  -0.41%  '

$ws.Range("D38").Value = '2.09'
$ws.Range("E38").Value = '  -3.20%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.20'
$ws.Range("E39").Value = '  +1.98%  '

$ws.Range("D40").Value = '50.35'
$ws.Range("E40").Value = '  +0.11%  '

$ws.Range("E41").Value = '  +8.55%  '

$ws.Range("E42").Value = '  -3.71%  '

$ws.Range("D43").Value = '400.32'
$ws.Range("E43").Value = '  -5.40%  '

$ws.Range("E44").Value = '  +3.45%  '

$ws.Range("D45").Value = '0.0353'
$ws.Range("E45").Value = '  -0.75%  '

$ws.Range("E46").Value = '  -5.71%  '

$ws.Range("D47").Value = '2.718.15'
$ws.Range("E47").Value = '  -2.89%  '

$ws.Range("D48").Value = '133.23'
$ws.Range("E48").Value = '  +2.54%  '

$ws.Range("E49").Value = '  +0.07%  '

$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").Value = '2.18'
$ws.Range("E50").Value = '  +1.47%  '

$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").Value = '0.108'
$ws.Range("E51").Value = '  -0.55%  '
